$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.442.54"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "2.610.19"
$ws.Range("E3").Value = "  +10.05%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.34%  "

$ws.Range("E7").Value = "  +5.94%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.95%  "

$ws.Range("D13").Value = "3.009.49"
$ws.Range("E13").Value = "  +10.08%  "

$ws.Range("E14").Value = "  +1.93%  "

$ws.Range("D15").Value = "2.608.64"
$ws.Range("E15").Value = "  +10.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.905"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +10.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.17%  "

$ws.Range("D18").Value = "46.600.27"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.73%  "

$ws.Range("E20").Value = "  +4.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +15.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +35.20%  "

$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("E28").Value = "  +7.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("E30").Value = "  +3.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.78%  "

$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +22.53%  "

$ws.Range("E34").Value = "  +5.26%  "

$ws.Range("E35").Value = "  +7.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.80%  "

$ws.Range("E37").Value = "  +3.99%  "

$ws.Range("E38").Value = "  +5.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.82%  "

$ws.Range("E40").Value = "  +6.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.37%  "

$ws.Range("E42").Value = "  +7.56%  "

$ws.Range("D43").Value = "2.050.72"
$ws.Range("E43").Value = "  +6.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +31.89%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.55%  "

$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.51%  "

$ws.Range("D51").Value = "2.865.35"
$ws.Range("E51").Value = "  +10.03%  "
